$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings) - temporarily force text format so
# Excel does not auto-convert these month/year-looking strings into date
# serials, then restore the original (default/Normal) cell style so the
# cell formatting is left exactly as it was before the edit.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "September 2024"
$ws.Range("A1").Style = "Normal"

$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "October 2024"
$ws.Range("G1").Style = "Normal"

# Update the numeric data row
$ws.Range("A2").Value = -0.004
$ws.Range("B2").Value = 0.077
$ws.Range("C2").Value = 0.1
$ws.Range("D2").Value = 0.068
$ws.Range("E2").Value = 0.022
$ws.Range("F2").Value = -0.104
$ws.Range("G2").Value = 0.155
